$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right after
#    the title heading (paragraph 2).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Insert a new bold "Play 243 Crystal Fruits Free - Unique Graphics &
#    Features" paragraph right before the final paragraph of the document
#    (the one that currently holds the italic image-prompt text).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertStart = $lastPara.Range.Start
$insertPoint = $d.Range($insertStart, $insertStart)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 243 Crystal Fruits Free - Unique Graphics &amp; Features</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insertPoint.InsertXML($newParaXml)

# InsertXML splitting at the start of the (still) final paragraph leaves a
# spare empty paragraph behind the newly inserted one - drop it.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -eq 1 -and $t -eq [char]13) {
        if ($i -gt 1 -and $d.Paragraphs.Item($i - 1).Range.Text -like "*Play 243 Crystal Fruits Free*") {
            $p.Range.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Swap the text of the (still) final paragraph - the image-generation
#    prompt - for the new meta-description sentence, keeping its italic
#    run formatting intact.
# ---------------------------------------------------------------------------
$oldText = 'Create a colorful cartoon-style image featuring a happy Maya warrior wearing glasses for the game "243 Crystal Fruits". The Maya warrior should hold a crystal fruit in one hand and a Wild symbol in the other hand. The background should be a vibrant jungle with cascading crystals falling from above. The image should convey a sense of fun and excitement, highlighting the unique and innovative graphic style of the game while also showcasing the ancient theme.'
$newText = 'Experience the innovative cascade system and multiplier on winnings with 243 Crystal Fruits. Enjoy the unique graphics and play for free now.'

[void]$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
